$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4: "Indicador" -> "KPI", drop the top border (card merges with header above) ---
$ws.Range("B4").Value2 = "KPI"
$ws.Range("B4").Borders.Item(8).LineStyle = -4142

# --- Row 5: "Indicador" -> "KPI", drop the top border ---
$ws.Range("B5").Value2 = "KPI"
$ws.Range("B5").Borders.Item(8).LineStyle = -4142

# --- Row 6: "Indicador" -> "MÉTRICA" (new label) ---
$ws.Range("B6").Value2 = "MÉTRICA"

# --- Row 10: "Indicador" -> "KPI" ---
$ws.Range("B10").Value2 = "KPI"

# --- Row heights ---
$ws.Rows.Item(8).RowHeight = 45
$ws.Rows.Item(11).RowHeight = 75

# --- Sheet view: reset scroll position (drop topLeftCell="F1") and move selection/active cell to B11 ---
$excel.ActiveWindow.ScrollColumn = 1
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("B11").Select()
